# Applies the COLLECTIONS.pptx "WHY COLLECTIONS" slide edit:
#   - removes the paragraph "Group of objects into a single entity."
#   - splits "To store similar or dissimilar types of objects into a
#     single unit." into two runs (break right before "unit.")
#   - splits "To represent group of individual objects as a single
#     entity." into two runs (break right before "represent")
#
# Locates the slide/shape by content instead of hard-coded indices so the
# script is resilient to ordering differences.

$p = $ppt.ActivePresentation

$targetShape = $null
$targetSlide = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $txt = $shape.TextFrame.TextRange.Text
            if ($txt -like "*To store similar or dissimilar types of objects into a single unit.*") {
                $targetShape = $shape
                $targetSlide = $slide
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange

# --- Step 1: remove the middle paragraph ("Group of objects into a single
#     entity.") entirely; this merges what was paragraph 3 up to become
#     paragraph 2. ---------------------------------------------------------
$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $paraText = $para.Text.TrimEnd([char]13, [char]10)
    if ($paraText -eq "Group of objects into a single entity.") {
        $para.Delete()
        break
    }
}

# --- Step 2: split "...into a single unit." into two runs, breaking right
#     before "unit." -------------------------------------------------------
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $paraText = $para.Text.TrimEnd([char]13, [char]10)
    if ($paraText -eq "To store similar or dissimilar types of objects into a single unit.") {
        $splitAt = $paraText.IndexOf("unit.")
        $tail = $para.Characters($splitAt + 1, $paraText.Length - $splitAt)
        $tail.Text = $tail.Text
        break
    }
}

# --- Step 3: split "To represent group of individual objects as a single
#     entity." into two runs, breaking right before "represent" -----------
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $paraText = $para.Text.TrimEnd([char]13, [char]10)
    if ($paraText -eq "To represent group of individual objects as a single entity.") {
        $splitAt = $paraText.IndexOf("represent")
        $tail = $para.Characters($splitAt + 1, $paraText.Length - $splitAt)
        $tail.Text = $tail.Text
        break
    }
}
